# Append 9 new report rows (205-213) to the "Relatórios" sheet — the
# user now asks the form to generate several reports at once (min 1,
# max 10) using the same field values, so the save route appends one
# row per generated report number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part Number (C) and Semana (E) are plain digit strings; left at the
# default "General" format Excel would silently re-type them as numbers
# on write, so force Text format first, same as the rest of the sheet
# (hence the sheet's numberStoredAsText warnings being suppressed).
$ws.Range("C205:C213").NumberFormat = "@"
$ws.Range("E205:E213").NumberFormat = "@"

$rows = @(
    @("C2026.0204", "27/01/2026", "53327555", "RINFORZO EST. IN CINTURA PORTA ANT. DX", "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0205", "27/01/2026", "53327555", "RINFORZO EST. IN CINTURA PORTA ANT. DX", "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0206", "27/01/2026", "53327555", "RINFORZO EST. IN CINTURA PORTA ANT. DX", "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0207", "27/01/2026", "53327555", "RINFORZO EST. IN CINTURA PORTA ANT. DX", "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0208", "27/01/2026", "53327555", "RINFORZO EST. IN CINTURA PORTA ANT. DX", "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0209", "27/01/2026", "53490369", "MANCAL DO LONGHERONE LT",                "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "ANÁLISE DIMENSIONAL", ""),
    @("C2026.0210", "27/01/2026", "53490369", "MANCAL DO LONGHERONE LT",                "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "ANÁLISE DIMENSIONAL", ""),
    @("C2026.0211", "27/01/2026", "51947034", "STAFFA COMPL FISS INF PARAFANGO RT",     "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", ""),
    @("C2026.0212", "27/01/2026", "51947034", "STAFFA COMPL FISS INF PARAFANGO RT",     "5", "QUALIDADE", "luis", "2º TURNO", "METRASCAN", "INSP LAYOUT", "")
)

$startRow = 205
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
    $ws.Range("F$r").Value = $data[5]
    $ws.Range("G$r").Value = $data[6]
    $ws.Range("H$r").Value = $data[7]
    $ws.Range("I$r").Value = $data[8]
    $ws.Range("J$r").Value = $data[9]
    $ws.Range("K$r").Value = $data[10]
    # "Selecionado" column is boolean FALSE for every newly generated report.
    $ws.Range("L$r").Value = $false
}
